# Update cryptocurrency price/volume figures (Price column D, Volume(1h) column E).
# Source cells are text-formatted (prices use "." as a thousands separator, not a
# decimal point), so force NumberFormat = "@" before writing to keep Excel from
# reinterpreting the strings as numbers and stripping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.005.19'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.630.68'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.25'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0619'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.61'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0788'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.858.26'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.642.05'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.76%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.012.53'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0₃0740'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.48'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.22%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.19'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.55'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.08'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.135'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.21'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.74'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.28'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.53%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.08%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.13'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.04%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.98%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.98%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.122.41'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.44'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.521'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0154'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.16'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.770'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -3.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.767.84'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.45%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0530'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.53'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.60%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.49'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.01%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.07%  '
